$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 6 (shifts old rows 6-9 down to 7-10,
# and shifts/updates the existing merged ranges automatically).
$ws.Rows.Item(6).Insert()

# --- Merge A5:A6 (cylinder=4 group now spans two rows, like the 6- and
#     8-cylinder groups below it). Do this before re-applying the
#     distinct per-cell formatting below, since merging harmonizes the
#     style of the whole range to a single one.
$ws.Range("A5:A6").MergeCells = $true

# --- Fix up styling for the brand-new row 6 by copying formats from
#     the analogous cells in row 5 / the H column "blank bordered" cell.
$ws.Range("H5").Copy()
$ws.Range("A6").PasteSpecial(-4122)
$ws.Range("H6").PasteSpecial(-4122)

$ws.Range("D5").Copy()
$ws.Range("D6").PasteSpecial(-4122)

$ws.Range("B5:C5").Copy()
$ws.Range("B6:C6").PasteSpecial(-4122)

$ws.Range("E5:G5").Copy()
$ws.Range("E6:G6").PasteSpecial(-4122)

$excel.CutCopyMode = 0

# --- Row 5 (cylinder = 4) updated summary values.
$ws.Range("D5").Value = 91
$ws.Range("E5").Value = ""
$ws.Range("F5").Value = 2.14
$ws.Range("G5").Value = ""

# --- Row 6 (new second sub-row for cylinder = 4).
$ws.Range("B6").Value = 1
$ws.Range("C6").Value = 1
$ws.Range("D6").Value = 81.8
$ws.Range("E6").Value = 21.87235698318771
$ws.Range("F6").Value = 2.3003
$ws.Range("G6").Value = 0.5982073312080948

# --- Row 7 (cylinder = 6, first sub-row) updated summary values.
$ws.Range("D7").Value = 131.6666666666667
$ws.Range("E7").Value = 37.52776749732568
$ws.Range("F7").Value = 2.755
$ws.Range("G7").Value = 0.1281600561797629

# --- Row 8 (cylinder = 6, second sub-row) updated summary values.
$ws.Range("D8").Value = 115.25
$ws.Range("E8").Value = 9.178779875342908
$ws.Range("F8").Value = 3.38875
$ws.Range("G8").Value = 0.1162163929916946

# --- Row 9 (cylinder = 8) updated summary values.
$ws.Range("D9").Value = 209.2142857142857
$ws.Range("E9").Value = 50.97688551827051
$ws.Range("F9").Value = 3.999214285714287
$ws.Range("G9").Value = 0.7594047444769265
